$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.847.61'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = '''1.893.85'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''0.7896'
$ws.Range('E5').Value = '  -5.17%  '
$ws.Range('D6').Value = '''243.85'
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '''0.3150'
$ws.Range('E8').Value = '  -4.12%  '
$ws.Range('D9').Value = '''25.30'
$ws.Range('E9').Value = '  -5.84%  '
$ws.Range('E10').Value = '  +2.29%  '
$ws.Range('D11').Value = '''0.08091'
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').Value = '''0.7665'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''5.536'
$ws.Range('E13').Value = '  +4.66%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '''1.933.88'
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('D15').Value = '''92.52'
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').Value = '''6.149'
$ws.Range('E16').Value = '  +4.19%  '
$ws.Range('D17').Value = '''29.875.89'
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('D18').Value = '''13.93'
$ws.Range('E18').Value = '  -1.89%  '
$ws.Range('D19').Value = '''243.96'
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('D20').Value = '''0.000007765'
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('D21').Value = '''2.167.17'
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').Value = '''8.154'
$ws.Range('E23').Value = '  +15.88%  '
$ws.Range('D24').Value = '''1.001'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '''0.1651'
$ws.Range('E25').Value = '  -5.36%  '
$ws.Range('D26').Value = '''9.402'
$ws.Range('E26').Value = '  +0.98%  '
$ws.Range('D27').Value = '''163.28'
$ws.Range('E27').Value = '  -1.82%  '
$ws.Range('D28').Value = '''18.74'
$ws.Range('E28').Value = '  -1.51%  '
$ws.Range('D29').Value = '''2.054'
$ws.Range('E29').Value = '  -2.36%  '
$ws.Range('D30').Value = '''1.403'
$ws.Range('E30').Value = '  +2.48%  '
$ws.Range('D31').Value = '''1.546'
$ws.Range('E31').Value = '  +1.55%  '
$ws.Range('D32').Value = '''4.492'
$ws.Range('E32').Value = '  +4.37%  '
$ws.Range('D33').Value = '''4.110'
$ws.Range('E33').Value = '  +0.51%  '
$ws.Range('D34').Value = '''0.05585'
$ws.Range('E34').Value = '  -7.56%  '
$ws.Range('D35').Value = '''1.270'
$ws.Range('E35').Value = '  -0.39%  '
$ws.Range('D36').Value = '''0.7421'
$ws.Range('E36').Value = '  +1.09%  '
$ws.Range('D37').Value = '''1.003'
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('D38').Value = '''2.618'
$ws.Range('E38').Value = '  -3.70%  '
$ws.Range('D39').Value = '''0.01927'
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('D40').Value = '''2.777'
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range('D41').Value = '''1.142.73'
$ws.Range('E41').Value = '  +13.44%  '
$ws.Range('D42').Value = '''74.10'
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('D43').Value = '''0.4432'
$ws.Range('E43').Value = '  -0.76%  '
$ws.Range('D44').Value = '''5.888'
$ws.Range('E44').Value = '  -1.03%  '
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('D46').Value = '''104.44'
$ws.Range('E46').Value = '  +1.95%  '
$ws.Range('D47').Value = '''1.000'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').Value = '''10.07'
$ws.Range('E48').Value = '  +2.22%  '
$ws.Range('D49').Value = '''1.880'
$ws.Range('E49').Value = '  -1.63%  '
$ws.Range('D50').Value = '''7.462'
$ws.Range('E50').Value = '  -1.74%  '
$ws.Range('E51').Value = '  +9.14%  '
